$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the cell to keep its text/general style while storing a
    # numeric-looking string as text (matches the source workbook, which
    # stores every Price/Volume cell as inline text, never as a number).
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '60.305.35'
$ws.Range('E2').Value = '  +2.18%  '

$ws.Range('D3').Value = '2.671.07'
$ws.Range('E3').Value = '  +1.51%  '

$ws.Range('E4').Value = '  +0.06%  '

Set-TextValue 'D5' '520.06'
$ws.Range('E5').Value = '  +1.16%  '

Set-TextValue 'D6' '145.64'
$ws.Range('E6').Value = '  +0.84%  '

$ws.Range('E7').Value = '  +0.22%  '

$ws.Range('E8').Value = '  +1.74%  '

$ws.Range('D9').Value = '2.688.59'
$ws.Range('E9').Value = '  +1.21%  '

$ws.Range('E10').Value = '  +1.53%  '

Set-TextValue 'D11' '0.105'
$ws.Range('E11').Value = '  -0.39%  '

$ws.Range('E12').Value = '  +0.72%  '

$ws.Range('E13').Value = '  +1.45%  '

$ws.Range('D14').Value = '3.143.94'
$ws.Range('E14').Value = '  +1.81%  '

$ws.Range('D15').Value = '60.311.27'
$ws.Range('E15').Value = '  +2.32%  '

Set-TextValue 'D16' '21.27'
$ws.Range('E16').Value = '  +0.77%  '

$ws.Range('D17').Value = '2.762.58'
$ws.Range('E17').Value = '  +4.38%  '

Set-TextValue 'D18' '0.0000138'
$ws.Range('E18').Value = '  +0.82%  '

Set-TextValue 'D19' '349.98'
$ws.Range('E19').Value = '  +1.55%  '

$ws.Range('E20').Value = '  -0.03%  '

Set-TextValue 'D21' '10.53'
$ws.Range('E21').Value = '  +1.42%  '

Set-TextValue 'D22' '6.30'
$ws.Range('E22').Value = '  +3.23%  '

Set-TextValue 'D23' '0.998'
$ws.Range('E23').Value = '  +0.03%  '

Set-TextValue 'D24' '62.70'
$ws.Range('E24').Value = '  +2.59%  '

Set-TextValue 'D25' '0.421'
$ws.Range('E25').Value = '  +0.04%  '

$ws.Range('E26').Value = '  +4.61%  '

$ws.Range('E27').Value = '  +0.39%  '

$ws.Range('D28').Value = '0.0₃0812'
$ws.Range('E28').Value = '  +0.74%  '

Set-TextValue 'D29' '7.24'
$ws.Range('E29').Value = '  +1.51%  '

Set-TextValue 'D30' '6.86'
$ws.Range('E30').Value = '  +6.42%  '

Set-TextValue 'D31' '0.999'
$ws.Range('E31').Value = '  +0.21%  '

Set-TextValue 'D32' '19.06'
$ws.Range('E32').Value = '  +0.84%  '

$ws.Range('E33').Value = '  +0.72%  '

Set-TextValue 'D34' '148.33'
$ws.Range('E34').Value = '  -1.18%  '

Set-TextValue 'D35' '4.30'
$ws.Range('E35').Value = '  +6.57%  '

Set-TextValue 'D36' '0.950'
$ws.Range('E36').Value = '  -6.98%  '

Set-TextValue 'D37' '1.22'
$ws.Range('E37').Value = '  +5.80%  '

$ws.Range('E38').Value = '  +10.44%  '

Set-TextValue 'D39' '0.870'
$ws.Range('E39').Value = '  +1.74%  '

Set-TextValue 'D40' '36.66'
$ws.Range('E40').Value = '  +0.44%  '

Set-TextValue 'D41' '3.69'
$ws.Range('E41').Value = '  -0.27%  '

Set-TextValue 'D42' '280.64'
$ws.Range('E42').Value = '  -0.30%  '

Set-TextValue 'D43' '0.0990'
$ws.Range('E43').Value = '  +0.38%  '

Set-TextValue 'D46' '0.609'
$ws.Range('E46').Value = '  -0.76%  '

$ws.Range('D47').Value = '2.124.75'
$ws.Range('E47').Value = '  +7.23%  '

Set-TextValue 'D48' '0.0541'
$ws.Range('E48').Value = '  +0.72%  '

Set-TextValue 'D49' '4.84'
$ws.Range('E49').Value = '  +4.56%  '

Set-TextValue 'D50' '0.0234'
$ws.Range('E50').Value = '  +1.95%  '

Set-TextValue 'D51' '10.45'
$ws.Range('E51').Value = '  +1.68%  '

# Rows 44/45 swap order: row44 becomes FirstDigitalUSD, row45 becomes EnergySwap
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D44' '0.997'
$ws.Range('E44').Value = '  +0.29%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D45' '19.96'
$ws.Range('E45').Value = '  +2.15%  '
